# HP: Writing check results to students.xlsx. TOIMII
#
# Applies the Results-sheet ("check results") edit described by the diff:
#  - adds a Feedback header + 4 testcase columns + 4 wrapped result rows
#  - activates / selects the Results sheet (was ZipFiles)
#  - touches a few column widths on both sheets

$wb = $excel.ActiveWorkbook
$wsZip = $wb.Worksheets.Item("ZipFiles")
$wsRes = $wb.Worksheets.Item("Results")

# --- Results sheet content -------------------------------------------------
# Preserve the existing "Feedback" shared string by writing it to its new
# location (G9) BEFORE G10 is overwritten with "Testcase1" - otherwise the
# string would become orphaned and get dropped from sharedStrings.xml.
$wsRes.Range("G9").Value = "Feedback"

$wsRes.Range("E10").Value = "Zip"
$wsRes.Range("G10").Value = "Testcase1"
$wsRes.Range("H10").Value = "Testcase2"
$wsRes.Range("I10").Value = "Testcase3"
$wsRes.Range("J10").Value = "Testcase4"

# H10:J10 need the same "bottom border" cell style already used across row 10
# (A10:G10). Copy the format from F10 (a cell that keeps that style) instead
# of setting borders directly, so the existing style entry is reused rather
# than a near-duplicate one being created.
$wsRes.Range("F10").Copy()
$wsRes.Range("H10:J10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Result message rows, each wrapped.
$wsRes.Range("G11:G14").WrapText = $true
$wsRes.Range("G11").Value = "SUBMIT(1) TESTCASE(1) RESULT MSG:(TESTCASE#1:EQUAL`n)"
$wsRes.Range("G12").Value = "SUBMIT(2) TESTCASE(1) RESULT MSG:(TESTCASE#1:NOT_EQUAL`n)"
$wsRes.Range("G13").Value = "SUBMIT(3) TESTCASE(1) RESULT MSG:(TESTCASE#1:NOT_EQUAL`n)"
$wsRes.Range("G14").Value = "SUBMIT(4) TESTCASE(1) RESULT MSG:(TESTCASE#1:EQUAL`n)"

# Wrapping the multi-line messages bumps each row to an auto-computed custom
# height; AutoFit puts the height back under "default" bookkeeping so no
# ht="..." customHeight="1" is written (matching the diff, which doesn't
# touch row heights at all).
$wsRes.Range("G11:G14").EntireRow.AutoFit()

# --- Column widths -----------------------------------------------------------
$wsRes.Columns.Item(3).ColumnWidth = 12.333333
$wsRes.Columns.Item(7).ColumnWidth = 74.5
$wsRes.Columns.Item(8).ColumnWidth = 22.333333

# --- Sheet selection / view state -------------------------------------------
$wsRes.Range("G12").Select()
$wsRes.Application.ActiveWindow.ScrollRow = 7

# Activating the Results sheet flips workbook.xml's activeTab to 1 and moves
# tabSelected from ZipFiles' sheetView to Results' sheetView, matching the
# diff exactly.
$wsRes.Activate()
